# Backup QR Scanner data - 14/08/2025, 8:24:57 AM
#
# The scanner log sheet originally held two rows of scan events.
# The second (newer) scan entry's data is merged into row 2, and the
# now-redundant row 3 is removed. The sheet is also renamed to match
# the subject being tracked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: "Scanner" -> "Anatomy"
$ws.Name = "Anatomy"

# Update row 2 with the values that used to live in row 3
# (keep these as text, matching the original t="str" storage)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "373739"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "08:23:46"

# Remove the now-duplicate third row entirely (shrinks dimension to A1:F2)
$ws.Rows(3).Delete()
